$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 342 all hold the same date serial
# number (45203 -> 2023-10-04). This update bumps the value by one day
# to 45204 (2023-10-05) for every one of those rows.
$ws.Range("C2:C342").Value = 45204
